$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.771.77'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '''1.597.12'
$ws.Range('E3').Value = '  -1.70%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''209.20'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('E6').Value = '  -2.06%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '''22.42'
$ws.Range('E8').Value = '  -2.62%  '
$ws.Range('D9').Value = '''0.254'
$ws.Range('E9').Value = '  -1.53%  '
$ws.Range('E10').Value = '  -1.74%  '
$ws.Range('D11').Value = '''0.0867'
$ws.Range('E11').Value = '  -1.71%  '
$ws.Range('D12').Value = '''1.823.79'
$ws.Range('E12').Value = '  -1.72%  '
$ws.Range('D13').Value = '''1.596.57'
$ws.Range('E13').Value = '  -2.01%  '
$ws.Range('E14').Value = '  -2.61%  '
$ws.Range('D15').Value = '''0.533'
$ws.Range('E15').Value = '  -3.50%  '
$ws.Range('D16').Value = '''27.750.97'
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('D17').Value = '''63.55'
$ws.Range('E17').Value = '  -1.46%  '
$ws.Range('D18').Value = '''219.99'
$ws.Range('E18').Value = '  -3.13%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '''7.39'
$ws.Range('E19').Value = '  -2.73%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '''0.0₃0697'
$ws.Range('E20').Value = '  -2.53%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').Value = '''4.17'
$ws.Range('E22').Value = '  -3.56%  '
$ws.Range('D23').Value = '''9.82'
$ws.Range('E23').Value = '  -1.07%  '
$ws.Range('E24').Value = '  -3.82%  '
$ws.Range('D25').Value = '''154.08'
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').Value = '''7.18'
$ws.Range('E26').Value = '  +4.01%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').Value = '''15.21'
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('E29').Value = '  -4.28%  '
$ws.Range('E30').Value = '  -1.21%  '
$ws.Range('D31').Value = '''0.0474'
$ws.Range('E31').Value = '  -1.36%  '
$ws.Range('E32').Value = '  -4.39%  '
$ws.Range('D33').Value = '''1.377.80'
$ws.Range('E33').Value = '  -3.06%  '
$ws.Range('E34').Value = '  -3.15%  '
$ws.Range('E35').Value = '  -3.56%  '
$ws.Range('D36').Value = '''0.979'
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('E38').Value = '  -0.58%  '
$ws.Range('E39').Value = '  -2.67%  '
$ws.Range('D40').Value = '''0.831'
$ws.Range('E40').Value = '  -1.98%  '
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('E42').Value = '  -3.11%  '
$ws.Range('D43').Value = '''64.66'
$ws.Range('E43').Value = '  -0.75%  '
$ws.Range('D44').Value = '''2.18'
$ws.Range('E44').Value = '  +2.54%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '''5.26'
$ws.Range('E45').Value = '  -2.10%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '''1.75'
$ws.Range('E46').Value = '  -1.55%  '
$ws.Range('D47').Value = '''1.734.06'
$ws.Range('E47').Value = '  -1.72%  '
$ws.Range('D48').Value = '''86.79'
$ws.Range('E48').Value = '  -2.71%  '
$ws.Range('D49').Value = '''0.0₆0101'
$ws.Range('E49').Value = '  -1.30%  '
$ws.Range('E50').Value = '  -3.25%  '
$ws.Range('D51').Value = '''0.0497'
$ws.Range('E51').Value = '  -1.00%  '
